# Actualización automática 2025-09-30 15:30:09
#
# Updates the "PORCELANATO" sale for ALTAMIRANO MARCATOMA EDISON PAULINO and
# the "PIEDRA SINTERIZADA" sale for ARMIJOS BUSTAMANTE FRANCISCO RAFAEL (both
# asesor ALMEIDA CUATIN JHONATHANN CARLOS) as well as the septiembre figure
# for SANTANA JIMENEZ MARIA ELENA, then refreshes every row/column total and
# the derived "CUMPLIMIENTO MENSUAL" sheet that depends on them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (sales by product group)
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M2").Value = 6231.33     # PORCELANATO - ALTAMIRANO MARCATOMA EDISON PAULINO
$wsGrupo.Range("L3").Value = 537.34      # PIEDRA SINTERIZADA - ARMIJOS BUSTAMANTE FRANCISCO RAFAEL
$wsGrupo.Range("D28").Value = 1831.68    # 240X80 PORCELANATO - SANTANA JIMENEZ MARIA ELENA

# Footer counter of non-zero rows per column: PIEDRA SINTERIZADA gained one.
$wsGrupo.Range("L35").Value = "3 de 33"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (monthly sales) - septiembre column mirrors the above
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F2").Value = 6231.33   # ALTAMIRANO MARCATOMA EDISON PAULINO
$wsMensual.Range("F3").Value = 687.03    # ARMIJOS BUSTAMANTE FRANCISCO RAFAEL
$wsMensual.Range("F28").Value = 1831.68  # SANTANA JIMENEZ MARIA ELENA

# Column total for septiembre
$wsMensual.Range("F35").Value = 29705.29

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" (monthly compliance) - recompute VENTA /
# POR CUMPLIR / CUMPLIMIENTO for the affected groups plus the TOTAL row
# ---------------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3: 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 3592.51
$wsCumpl.Range("E3").Value = 5242.06354940916
$wsCumpl.Range("F3").Value = 0.4066421519848302

# Row 11: PIEDRA SINTERIZADA
$wsCumpl.Range("D11").Value = 4298.72
$wsCumpl.Range("E11").Value = -1376.49541814726
$wsCumpl.Range("F11").Value = 1.471043679084562

# Row 12: PORCELANATO
$wsCumpl.Range("D12").Value = 20190.39
$wsCumpl.Range("E12").Value = 2243.3653751766
$wsCumpl.Range("F12").Value = 0.9000004529933081

# Row 15: TOTAL
$wsCumpl.Range("D15").Value = 29962.52
$wsCumpl.Range("E15").Value = 8780.498813395923
$wsCumpl.Range("F15").Value = 0.7733656518691323
